$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'31.088.08"
$ws.Range("E2").Value = "  +1.28%  "
$ws.Range("D3").Value = "'1.957.05"
$ws.Range("E3").Value = "  +0.48%  "
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").Value = "'246.47"
$ws.Range("E5").Value = "  -0.19%  "
$ws.Range("D6").Value = "'1.002"
$ws.Range("E6").Value = "  +0.16%  "
$ws.Range("D7").Value = "'0.4921"
$ws.Range("E7").Value = "  +1.66%  "
$ws.Range("D8").Value = "'0.2978"
$ws.Range("E8").Value = "  +1.32%  "
$ws.Range("D9").Value = "'0.06850"
$ws.Range("E9").Value = "  +0.43%  "
$ws.Range("D10").Value = "'19.22"
$ws.Range("E10").Value = "  -1.01%  "
$ws.Range("D11").Value = "'107.96"
$ws.Range("E11").Value = "  -4.12%  "
$ws.Range("D12").Value = "'1.946.54"
$ws.Range("D13").Value = "'0.07765"
$ws.Range("E13").Value = "  +1.19%  "
$ws.Range("D14").Value = "'5.470"
$ws.Range("E14").Value = "  -0.86%  "
$ws.Range("D15").Value = "'0.7087"
$ws.Range("E15").Value = "  +2.73%  "
$ws.Range("D16").Value = "'282.23"
$ws.Range("E16").Value = "  -4.80%  "
$ws.Range("D17").Value = "'31.101.31"
$ws.Range("E17").Value = "  +1.19%  "
$ws.Range("D18").Value = "'13.32"
$ws.Range("E18").Value = "  +0.26%  "
$ws.Range("D19").Value = "'0.000007772"
$ws.Range("E19").Value = "  +0.97%  "
$ws.Range("D20").Value = "'2.201.58"
$ws.Range("E20").Value = "  -0.01%  "
$ws.Range("D21").Value = "'1.001"
$ws.Range("E21").Value = "  +0.14%  "
$ws.Range("D22").Value = "'5.508"
$ws.Range("E22").Value = "  -2.96%  "
$ws.Range("E23").Value = "  +0.25%  "
$ws.Range("D24").Value = "'6.510"
$ws.Range("E24").Value = "  -1.41%  "
$ws.Range("D25").Value = "'9.839"
$ws.Range("E25").Value = "  -0.32%  "
$ws.Range("D26").Value = "'169.82"
$ws.Range("E26").Value = "  +0.89%  "
$ws.Range("D27").Value = "'20.07"
$ws.Range("E27").Value = "  -1.34%  "
$ws.Range("D28").Value = "'2.214"
$ws.Range("E28").Value = "  +0.92%  "
$ws.Range("D29").Value = "'0.1056"
$ws.Range("E29").Value = "  -2.85%  "
$ws.Range("D30").Value = "'1.431"
$ws.Range("E30").Value = "  -0.41%  "
$ws.Range("D31").Value = "'1.586"
$ws.Range("E31").Value = "  -0.38%  "
$ws.Range("D32").Value = "'4.583"
$ws.Range("E32").Value = "  -2.77%  "
$ws.Range("D33").Value = "'4.459"
$ws.Range("E33").Value = "  -1.39%  "
$ws.Range("D34").Value = "'0.04982"
$ws.Range("E34").Value = "  -2.18%  "
$ws.Range("D35").Value = "'0.7596"
$ws.Range("E35").Value = "  -2.32%  "
$ws.Range("D36").Value = "'1.186"
$ws.Range("E36").Value = "  +2.11%  "
$ws.Range("D37").Value = "'2.734"
$ws.Range("E37").Value = "  +0.05%  "
$ws.Range("D38").Value = "'0.02041"
$ws.Range("E38").Value = "  -2.18%  "
$ws.Range("D39").Value = "'2.707"
$ws.Range("E39").Value = "  +0.20%  "
$ws.Range("D40").Value = "'2.178"
$ws.Range("E40").Value = "  +5.86%  "
$ws.Range("D41").Value = "'6.504"
$ws.Range("E41").Value = "  +9.85%  "
$ws.Range("D42").Value = "'74.69"
$ws.Range("E42").Value = "  +6.37%  "
$ws.Range("D43").Value = "'0.4512"
$ws.Range("E43").Value = "  +1.01%  "
$ws.Range("D46").Value = "'8.172"
$ws.Range("E46").Value = "  +10.73%  "
$ws.Range("D47").Value = "'1.002"
$ws.Range("E47").Value = "  -0.10%  "
$ws.Range("D48").Value = "'977.24"
$ws.Range("D49").Value = "'9.417"
$ws.Range("E49").Value = "  +0.12%  "
$ws.Range("E50").Value = "  +0.83%  "
$ws.Range("D51").Value = "'0.2590"
$ws.Range("E51").Value = "  +2.35%  "

# Rows 44 and 45 swap order (Quant <-> TrustWalletToken) with refreshed values
$ws.Range("B44").Value = "TrustWalletToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D44").Value = "'0.8860"
$ws.Range("E44").Value = "  +1.41%  "

$ws.Range("B45").Value = "Quant"
$ws.Range("C45").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D45").Value = "'109.43"
$ws.Range("E45").Value = "  -1.87%  "
